# Update the table style ("Data Sources from LFX" tables) from the old
# custom style {A1FA7028-C1A7-44A9-B84D-1D3E6534E97C} to the new style
# {AA5A75A7-0268-4E23-BE90-BAA816196550} on every slide that contains a
# table using it.

$p = $ppt.ActivePresentation

$oldStyleId = "{A1FA7028-C1A7-44A9-B84D-1D3E6534E97C}"
$newStyleId = "{AA5A75A7-0268-4E23-BE90-BAA816196550}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
